$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price/volume data to reflect latest values scraped on Mon Nov 11 18:33:37 UTC 2024

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "84.760.87"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  +4.83%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.292.56"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  +1.97%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.01%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "217.64"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.64%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "632.82"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.58%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.317"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +13.18%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.998"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.03%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.590"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -3.09%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "3.287.53"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  +1.95%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.593"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -4.70%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0000273"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -1.28%  "

# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -0.08%  "

# Row 14
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  +1.94%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "34.20"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  +3.56%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "5.40"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.38%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "84.538.70"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  +4.87%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.282.07"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  +1.90%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "14.54"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.03%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "3.15"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +2.28%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "437.48"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -2.78%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.14"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -2.88%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "5.24"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -3.46%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "7.39"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +2.89%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "5.49"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +12.43%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "12.14"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +9.11%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "3.459.94"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  +1.95%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "77.65"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.49%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.0000130"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  +0.72%  "

# Row 30
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -0.31%  "

# Row 31
$ws.Range("B31").Value = "Bittensor"
$ws.Range("C31").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "601.31"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +5.26%  "

# Row 32
$ws.Range("B32").Value = "InternetComputer(DFINITY)"
$ws.Range("C32").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.18"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -2.76%  "

# Row 33
$ws.Range("B33").Value = "Cronos"
$ws.Range("C33").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.160"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +28.46%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.00"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  +0.25%  "

# Row 35
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  +2.47%  "

# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.03"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -1.04%  "

# Row 37
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.02%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "23.12"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -3.18%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "6.39"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +8.75%  "

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.14%  "

# Row 41
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -1.02%  "

# Row 42
$ws.Range("B42").Value = "WhiteBITCoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "20.96"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +2.89%  "

# Row 43
$ws.Range("B43").Value = "Stacks"
$ws.Range("C43").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.04"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  +9.97%  "

# Row 44
$ws.Range("B44").Value = "dogwifhat"
$ws.Range("C44").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.04"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  +9.73%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "159.92"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -2.75%  "

# Row 46
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +0.02%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "188.42"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -2.44%  "

# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.84%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "44.90"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  +2.96%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.783"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.77%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "26.41"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +1.89%  "
